# GanttChart.xlsx edit — week 5 update:
#  - Daily report: add RT60/EDT analysis note for day 5 (row 6), taller row
#  - 10 Week: add "Extended" / "Cancelled/Delayed" legend entries + swatches,
#    recolor several Gantt bars from Todo(blue) to Done(green),
#    add an Extended(orange) bar + note, and a Cancelled/Delayed(red) bar + note
#  - Active sheet switches from "Daily report" back to "10 Week"

$wb = $excel.ActiveWorkbook

$wsWeek  = $wb.Worksheets.Item("10 Week")
$wsDaily = $wb.Worksheets.Item("Daily report")

# ---------------------------------------------------------------------------
# Daily report sheet: log entry for day 5 in a new column C, and taller row
# ---------------------------------------------------------------------------
$wsDaily.Range("C6").Value = "Generated all audio. Organized RIR analysis folders. Modified MATLAB scripts for compatibility. Investigated and troubleshooted EDT/RT60 analysis issues. Experimented with different audio settings (sampling rate, sine sweep type, probe number, audio delay, volume). Identified persistent problems with audio analysis results."
$wsDaily.Range("C6").Copy()
$wsDaily.Range("B6").PasteSpecial(-4122)
$wsDaily.Application.CutCopyMode = 0
$wsDaily.Rows.Item(6).RowHeight = 101.5

$wsDaily.Range("C5").Select()

# ---------------------------------------------------------------------------
# 10 Week sheet: extend the Legends row with "Extended" (orange) and
# "Cancelled/Delayed" (red) swatches, recolor a few bars, and add notes
# ---------------------------------------------------------------------------

# Recolor several bars from "Todo" (blue) to "Done" (green) by copying the
# format from an existing green-filled cell.
$wsWeek.Range("Q1").Copy()
$wsWeek.Range("B3").PasteSpecial(-4122)
$wsWeek.Range("C4").PasteSpecial(-4122)
$wsWeek.Range("D4").PasteSpecial(-4122)
$wsWeek.Range("D5").PasteSpecial(-4122)
$wsWeek.Range("E5").PasteSpecial(-4122)
$wsWeek.Application.CutCopyMode = 0

# Note that goes with the newly-extended bar
$wsWeek.Range("N5").Value = "Gonna go longer than week 4, hopefully finish analysis with pretty graphs etc by end of Week 5"

$wsWeek.Range("R1").Value = "Extended"
$wsWeek.Range("T1").Value = "Cancelled/Delayed"

# Note that goes with the newly-cancelled/delayed bar
$wsWeek.Range("N6").Value = "This is now optional, SSC Mona is using stereo (some problem/setback?)"

# Orange swatch (theme Accent2) — set once, then propagate via format copy so
# the style table doesn't accumulate duplicate fills.
$wsWeek.Range("S1").Interior.ThemeColor = 6
$wsWeek.Range("S1").Copy()
$wsWeek.Range("F5").PasteSpecial(-4122)
$wsWeek.Application.CutCopyMode = 0

# Red swatch (solid red) — same propagation trick.
$wsWeek.Range("V1").Interior.Color = 255
$wsWeek.Range("V1").Copy()
$wsWeek.Range("F6").PasteSpecial(-4122)
$wsWeek.Application.CutCopyMode = 0

$wsWeek.Range("N6").Select()

# ---------------------------------------------------------------------------
# Make "10 Week" the active sheet/tab again
# ---------------------------------------------------------------------------
$wsWeek.Activate()
